$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Manatsa): reorder MARKS values and switch COURSE separator/order
$ws.Range("J2").Value = "[[86.75, 37.0, 90.45]]"
$ws.Range("K2").Value = "[{Comp Science : 2012} ,{Maths : 2010} ]"

# Row 3 (Grace): reorder MARKS values and switch COURSE separator/order
$ws.Range("J3").Value = "[[57.0, 56.05, 70.0]]"
$ws.Range("K3").Value = "[{Divorce : 2023} ,{Abuse : 2024} ]"
